$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-24 Saturday" "2026-01-25 Sunday"

Replace-Text "41×45=1845" "12×94=1128"
Replace-Text "23×64=1472" "24×28=672"
Replace-Text "45×58=2610" "47×87=4089"
Replace-Text "81×17=1377" "54×17=918"
Replace-Text "63×58=3654" "35×53=1855"

Replace-Text "47×81=3807" "99×64=6336"
Replace-Text "52×96=4992" "80×59=4720"
Replace-Text "39×65=2535" "76×55=4180"
Replace-Text "13×32=416" "63×94=5922"
Replace-Text "68×76=5168" "42×27=1134"

Replace-Text "34×59=2006" "12×71=852"
Replace-Text "42×52=2184" "33×42=1386"
Replace-Text "16×94=1504" "53×67=3551"
Replace-Text "37×79=2923" "66×55=3630"
Replace-Text "13×98=1274" "31×38=1178"

Replace-Text "25×68=1700" "94×79=7426"
Replace-Text "71×54=3834" "57×34=1938"
Replace-Text "80×81=6480" "78×39=3042"
Replace-Text "27×59=1593" "19×42=798"
Replace-Text "58×70=4060" "73×44=3212"

Replace-Text "35×14=490" "15×22=330"
Replace-Text "36×37=1332" "29×59=1711"
Replace-Text "69×39=2691" "44×36=1584"
Replace-Text "96×69=6624" "37×41=1517"
Replace-Text "63×56=3528" "93×29=2697"
